$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholder fields: "15/02/2023" -> "16/02/2023"
#    These live on the Slide Master and on every Slide Layout (12 places
#    total: 1 master + 11 layouts).
# ---------------------------------------------------------------------------
function Update-DateShapes($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Name -like "*Date Placeholder*") {
            if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
                if ($shp.TextFrame.TextRange.Text -eq "15/02/2023") {
                    $shp.TextFrame.TextRange.Text = "16/02/2023"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes $master.Shapes

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DateShapes $layouts.Item($i).Shapes
}

# ---------------------------------------------------------------------------
# 2) "About Ricky" slide, requirements paragraph: split the sentence
#    "...because their design is clean and intuitive." into three runs:
#      "...because "  +  "their designs are "  +  "clean and intuitive."
# ---------------------------------------------------------------------------
$targetShape = $null
$targetSlide = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -like "*Ricky requires a site*") {
                $targetSlide = $slide
                $targetShape = $shp
            }
        }
    }
}

if ($targetShape -ne $null) {
    $tr = $targetShape.TextFrame.TextRange
    $full = $tr.Text
    $anchorIdx = $full.IndexOf("Ricky requires a site")

    $oldTail = "their design is clean and intuitive."
    $tailPos0 = $full.IndexOf($oldTail, $anchorIdx)
    if ($tailPos0 -ge 0) {
        $tailPos1 = $tailPos0 + 1
        $sub = $tr.Characters($tailPos1, $oldTail.Length)
        $sub.Text = "their designs are clean and intuitive."

        $full2 = $tr.Text
        $cleanText = "clean and intuitive."
        $cleanPos0 = $full2.IndexOf($cleanText, $anchorIdx)
        $cleanPos1 = $cleanPos0 + 1
        $sub2 = $tr.Characters($cleanPos1, $cleanText.Length)
        $sub2.Text = "clean and intuitive."
    }
}
